$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newHeaders = @(
    "nGroup",
    "nCollaboration",
    "nLaneSet",
    "nLane",
    "nDataObject",
    "nDataObjectReference",
    "nDataStore",
    "nDataStoreReference",
    "nDataInput",
    "nDataOutput",
    "nExclusiveGateway",
    "nParallelGateway",
    "nInclusiveGateway",
    "nEventBasedGateway"
)

$startCol = 13  # M

# Copy the header-row formatting (bold font, boxed border, centered
# alignment) from the last existing header cell (L1) onto the new header
# cells, instead of re-deriving it property by property.
$ws.Range("L1").Copy()

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $startCol + $i
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $newHeaders[$i]
    $cell.PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item(2, $col).Value = 0
}
